$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2865
$ws.Range("C2").Value = "2022-11-13 17:16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1964"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "2022-11-13 16:21"

# Row 3
$ws.Range("B3").Value = 2274
$ws.Range("C3").Value = "2022-11-13 17:16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4542"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "2022-11-13 16:21"
